$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 48587
$ws.Range("C2").Value = 48587
$ws.Range("D2").Value = 48587
$ws.Range("E2").Value = 48587
$ws.Range("F2").Value = 48587
$ws.Range("G2").Value = 48587
$ws.Range("H2").Value = 48587
$ws.Range("I2").Value = 48587
$ws.Range("J2").Value = 48587
$ws.Range("K2").Value = 48587
$ws.Range("L2").Value = 48587
$ws.Range("M2").Value = 48587
$ws.Range("N2").Value = 48587
$ws.Range("O2").Value = 48587
$ws.Range("P2").Value = 48587
$ws.Range("Q2").Value = 48587
$ws.Range("B3").Value = 18482.26115607055
$ws.Range("C3").Value = [double]"-2.620644580818119e-16"
$ws.Range("D3").Value = [double]"3.743777972597313e-17"
$ws.Range("E3").Value = [double]"8.401564239285767e-17"
$ws.Range("F3").Value = [double]"3.743777972597313e-17"
$ws.Range("G3").Value = [double]"1.123133391779194e-16"
$ws.Range("H3").Value = [double]"2.293064008215854e-16"
$ws.Range("I3").Value = [double]"-8.490771448789062e-16"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = [double]"1.871888986298657e-17"
$ws.Range("L3").Value = [double]"-7.019583698619961e-18"
$ws.Range("M3").Value = [double]"2.807833479447985e-17"
$ws.Range("N3").Value = [double]"-5.615666958895969e-17"
$ws.Range("O3").Value = [double]"2.246266783558388e-16"
$ws.Range("P3").Value = [double]"-1.179290061368154e-15"
$ws.Range("Q3").Value = [double]"9.608055187486073e-17"
$ws.Range("B4").Value = 28191.71108331252
$ws.Range("C4").Value = 1.000010290977386
$ws.Range("D4").Value = 1.000010290977386
$ws.Range("E4").Value = 1.000010290977386
$ws.Range("F4").Value = 1.000010290977386
$ws.Range("G4").Value = 1.000010290977386
$ws.Range("H4").Value = 1.000010290977386
$ws.Range("I4").Value = 1.000010290977386
$ws.Range("J4").Value = 1.000010290977386
$ws.Range("K4").Value = 1.000010290977386
$ws.Range("L4").Value = 1.000010290977386
$ws.Range("M4").Value = 1.000010290977386
$ws.Range("N4").Value = 1.000010290977386
$ws.Range("O4").Value = 1.000010290977386
$ws.Range("P4").Value = 1.000010290977386
$ws.Range("Q4").Value = 1.000010290977386
$ws.Range("C5").Value = -1.378016816427322
$ws.Range("D5").Value = -1.422109209625639
$ws.Range("E5").Value = -0.2760353945003954
$ws.Range("F5").Value = -0.7653011471062969
$ws.Range("G5").Value = -1.277366054068879
$ws.Range("H5").Value = -2.598882644753269
$ws.Range("I5").Value = -1.662385164068257
$ws.Range("J5").Value = -0.4369634059966164
$ws.Range("K5").Value = -0.1613729815367844
$ws.Range("L5").Value = -0.08870602236863609
$ws.Range("M5").Value = -0.2998120537358804
$ws.Range("N5").Value = -0.470598222602185
$ws.Range("O5").Value = -1.999443984758019
$ws.Range("P5").Value = -2.968119362434912
$ws.Range("Q5").Value = -1.676597354793999
$ws.Range("B6").Value = 2207.065
$ws.Range("C6").Value = -0.6697668422407405
$ws.Range("D6").Value = -0.8160386215671275
$ws.Range("E6").Value = -0.2760353945003954
$ws.Range("F6").Value = -0.7653011471062969
$ws.Range("G6").Value = -1.243898881627937
$ws.Range("H6").Value = -0.7566869081585871
$ws.Range("I6").Value = -1.083600848764242
$ws.Range("J6").Value = -0.4369634059966164
$ws.Range("K6").Value = -0.1613729815367844
$ws.Range("L6").Value = -0.08870602236863609
$ws.Range("M6").Value = -0.2998120537358804
$ws.Range("N6").Value = -0.470598222602185
$ws.Range("O6").Value = 0.3477365382802091
$ws.Range("P6").Value = -0.4447741502816497
$ws.Range("Q6").Value = -0.7537898138438314
$ws.Range("B7").Value = 7358.74
$ws.Range("C7").Value = 0.03848313194584081
$ws.Range("D7").Value = -0.2436386217340888
$ws.Range("E7").Value = -0.2760353945003954
$ws.Range("F7").Value = -0.7653011471062969
$ws.Range("G7").Value = 0.2708729443511607
$ws.Range("H7").Value = 0.1135363685758193
$ws.Range("I7").Value = 0.1708255652211969
$ws.Range("J7").Value = -0.4369634059966164
$ws.Range("K7").Value = -0.1613729815367844
$ws.Range("L7").Value = -0.08870602236863609
$ws.Range("M7").Value = -0.2998120537358804
$ws.Range("N7").Value = -0.470598222602185
$ws.Range("O7").Value = 0.4530422998306055
$ws.Range("P7").Value = 0.1875123511544549
$ws.Range("Q7").Value = -0.1385847865437198
$ws.Range("B8").Value = 21254.705
$ws.Range("C8").Value = 0.7467331061324221
$ws.Range("D8").Value = 0.9348319661574616
$ws.Range("E8").Value = -0.2760353945003954
$ws.Range("F8").Value = 1.306675161511426
$ws.Range("G8").Value = 0.9488856805881467
$ws.Range("H8").Value = 0.8820104621843569
$ws.Range("I8").Value = 0.888990593190261
$ws.Range("J8").Value = 0.05041924447897494
$ws.Range("K8").Value = -0.1600238430902572
$ws.Range("L8").Value = -0.08783525164800124
$ws.Range("M8").Value = -0.1607730796183683
$ws.Range("N8").Value = 0.1332632143664955
$ws.Range("O8").Value = 0.5944031225457114
$ws.Range("P8").Value = 0.8304336408410139
$ws.Range("Q8").Value = 0.7842227544064475
$ws.Range("B9").Value = 385051.04
$ws.Range("C9").Value = 1.454983080319003
$ws.Range("D9").Value = 1.877608436470702
$ws.Range("E9").Value = 3.622723824275974
$ws.Range("F9").Value = 1.306675161511426
$ws.Range("G9").Value = 0.9961678677951344
$ws.Range("H9").Value = 1.718763612890517
$ws.Range("I9").Value = 1.628417085803147
$ws.Range("J9").Value = 12.65681544648191
$ws.Range("K9").Value = 16.17731873026641
$ws.Range("L9").Value = 15.89594578645218
$ws.Range("M9").Value = 12.17382441294384
$ws.Range("N9").Value = 10.03768053957526
$ws.Range("O9").Value = 0.7929365993549454
$ws.Range("P9").Value = 1.618374770306495
$ws.Range("Q9").Value = 1.707030295356615
